{"js": "// Template changes for Release 1.0\n// 1) Bump the AllHaem panel version referenced in the Test Methodology\n//    paragraph from \"v1\" to \"v3\".\n// 2) Update the cached DATE field result (\"Reported\" date) from\n//    1-Nov-2023 to 15-Nov-2023.\n\nconst body = context.document.body;\n\n// --- 1) \" v1) and sequenced\" -> \" v3) and sequenced\" ------------------\nconst versionHits = body.search(\" v1) and sequenced\", { matchCase: true });\nversionHits.load(\"items\");\nawait context.sync();\n\nfor (const hit of versionHits.items) {\n  hit.insertText(\" v3) and sequenced\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 2) \"1-Nov-2023\" -> \"15-Nov-2023\" (cached DATE field result) ------\nconst dateHits = body.search(\"1-Nov-2023\", { matchCase: true });\ndateHits.load(\"items\");\nawait context.sync();\n\nfor (const hit of dateHits.items) {\n  hit.insertText(\"15-Nov-2023\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Template changes for Release 1.0\n# 1) Bump the AllHaem panel version referenced in the Test Methodology\n#    paragraph from \"v1\" to \"v3\".\n# 2) Update the cached DATE field result (\"Reported\" date) from\n#    1-Nov-2023 to 15-Nov-2023.\n\n$d = $word.ActiveDocument\n\n# --- 1) \" v1) and sequenced\" -> \" v3) and sequenced\" ------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \" v1) and sequenced\"\n$find.Replacement.Text = \" v3) and sequenced\"\n$find.Execute([ref]\" v1) and sequenced\", [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]\" v3) and sequenced\", [ref]2)\n\n# --- 2) \"1-Nov-2023\" -> \"15-Nov-2023\" (cached DATE field result) ------\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"1-Nov-2023\"\n$find2.Replacement.Text = \"15-Nov-2023\"\n$find2.Execute([ref]\"1-Nov-2023\", [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]\"15-Nov-2023\", [ref]2)\n"}
